# Update the NFL weekly schedule table with the newly-scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schedule data: Teams (matchup), Time, Location (stadium)
$teams = @(
    "Tampa Bay  @  Atlanta",
    "NY Jets  vs.  Minnesota",
    "Cleveland  @  Washington",
    "Indianapolis  @  Jacksonville",
    "Miami  @  New England",
    "Buffalo  @  Houston",
    "Baltimore  @  Cincinnati",
    "Carolina  @  Chicago",
    "Las Vegas  @  Denver",
    "Arizona  @  San Francisco",
    "Green Bay  @  LA Rams",
    "NY Giants  @  Seattle",
    "Dallas  @  Pittsburgh",
    "New Orleans  @  Kansas City"
)

$times = @(
    "8:15 PM",
    "9:30 AM",
    "1:00 PM",
    "1:00 PM",
    "1:00 PM",
    "1:00 PM",
    "1:00 PM",
    "1:00 PM",
    "4:05 PM",
    "4:05 PM",
    "4:25 PM",
    "4:25 PM",
    "8:20 PM",
    "8:15 PM"
)

$locations = @(
    "Mercedes-Benz Stadium",
    "Tottenham Hotspur Stadium",
    "Northwest Stadium",
    "EverBank Stadium",
    "Gillette Stadium",
    "NRG Stadium",
    "Paycor Stadium",
    "Soldier Field",
    "Empower Field at Mile High",
    "Levi's Stadium",
    "SoFi Stadium",
    "Lumen Field",
    "Acrisure Stadium",
    "GEHA Field at Arrowhead Stadium"
)

# Overwrite the first 14 data rows (rows 2-15) with the fresh schedule.
for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $teams[$i]
    $ws.Cells.Item($row, 2).Value = $times[$i]
    $ws.Cells.Item($row, 3).Value = $locations[$i]
}

# The new schedule only has 14 games (was 16), so remove the two leftover rows.
$ws.Rows("16:17").Delete() | Out-Null

# Column A's best-fit width shifts slightly because of the new (longer) matchup text.
$ws.Columns("A:A").ColumnWidth = 25.7

# Reflect where the cursor ended up when the sheet was last saved.
$ws.Range("H19").Select() | Out-Null
